$wb = $excel.ActiveWorkbook
$wsAdd = $wb.Worksheets.Item("Add Devices Loop A")
$wsDel = $wb.Worksheets.Item("Delete Devices Loop A")

# --- Update "Add Devices Loop A" test data (row 7 / row 8) ---
# Row 7 becomes the LI800 entry with freshly measured drop values
$wsAdd.Range("A7").Value = "LI800"
$wsAdd.Range("B7").Value = "Other"
$wsAdd.Range("C7").Value = "LI800 - 1"
$wsAdd.Range("E7").Value = 289
$wsAdd.Range("F7").Value = 0.48
$wsAdd.Range("G7").Value = 0.89

# F7 picks up the quote-prefix formatting already used by F8
$wsAdd.Range("F8").Copy()
$wsAdd.Range("F7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 8 becomes the 801 CH entry (keeps its previously measured drop values)
$wsAdd.Range("A8").Value = "801 CH"
$wsAdd.Range("B8").Value = "Detectors"
$wsAdd.Range("C8").Value = "801 CH - 3"

# --- Make "Add Devices Loop A" the active/selected sheet & cell ---
$wsAdd.Activate()
$wsAdd.Range("B7").Select()
